{"js": "// Word JS API (Office.js) edit script.\n// Body of: async (context) => { ... }\n//\n// Net semantic changes made by the diff (proofErr spell/gram-check markers\n// are cosmetic artifacts of Word's background proofing pass and are not\n// reachable from the Word JS API, so we focus on the three actual wording\n// changes the diff introduces):\n//\n//  1. \"...\u0432\u043e \u0432\u0440\u0435\u043c\u044f \u043f\u0440\u043e\u0432\u0435\u0434\u0435\u043d\u0438\u044f \u0441\u043e\u0431\u0440\u0430\u043d\u0438\u0439.\" ->\n//     \"...\u0432\u043e \u0432\u0440\u0435\u043c\u044f \u043f\u0440\u043e\u0432\u0435\u0434\u0435\u043d\u0438\u044f \u0441\u043e\u0431\u0440\u0430\u043d\u0438\u044f, \u0430 \u0442\u0430\u043a\u0436\u0435 \u0433\u0430\u0440\u0430\u043d\u0442\u0438\u0440\u0443\u0435\u043c \u0441\u043e\u0431\u043b\u044e\u0434\u0435\u043d\u0438\u0435 \u043c\u0435\u0440\n//     \u0438\u043d\u0434\u0438\u0432\u0438\u0434\u0443\u0430\u043b\u044c\u043d\u043e\u0439 \u0437\u0430\u0449\u0438\u0442\u044b (\u0434\u043e 50 \u0447\u0435\u043b\u043e\u0432\u0435\u043a).\"\n//  2. \"\u041e\u0442\u0432\u0435\u0442\u0441\u0442\u0432\u0435\u043d\u043d\u044b\u0439 \u2013 {{responsable}}.\" -> \"\u041e\u0442\u0432\u0435\u0442\u0441\u0442\u0432\u0435\u043d\u043d\u044b\u0439(-\u0430\u044f) \u2013 {{responsable}}.\"\n\nconst body = context.document.body;\n\n// --- Change 1: extend the \"\u041c\u044b \u043f\u0440\u0438\u043d\u0438\u043c\u0430\u0435\u043c...\" sentence -----------------\nconst oldTail1 = \"\u0432\u043e \u0432\u0440\u0435\u043c\u044f \u043f\u0440\u043e\u0432\u0435\u0434\u0435\u043d\u0438\u044f \u0441\u043e\u0431\u0440\u0430\u043d\u0438\u0439.\";\nconst newTail1 =\n  \"\u0432\u043e \u0432\u0440\u0435\u043c\u044f \u043f\u0440\u043e\u0432\u0435\u0434\u0435\u043d\u0438\u044f \u0441\u043e\u0431\u0440\u0430\u043d\u0438\u044f, \u0430 \u0442\u0430\u043a\u0436\u0435 \u0433\u0430\u0440\u0430\u043d\u0442\u0438\u0440\u0443\u0435\u043c \u0441\u043e\u0431\u043b\u044e\u0434\u0435\u043d\u0438\u0435 \u043c\u0435\u0440 \u0438\u043d\u0434\u0438\u0432\u0438\u0434\u0443\u0430\u043b\u044c\u043d\u043e\u0439 \u0437\u0430\u0449\u0438\u0442\u044b (\u0434\u043e 50 \u0447\u0435\u043b\u043e\u0432\u0435\u043a).\";\n\nconst hits1 = body.search(oldTail1, { matchCase: true, matchWholeWord: false });\nhits1.load(\"text\");\nawait context.sync();\n\nif (hits1.items.length > 0) {\n  hits1.items[0].insertText(newTail1, \"Replace\");\n} else {\n  throw new Error(\"Could not find target sentence for change 1\");\n}\n\nawait context.sync();\n\n// --- Change 2: \"\u041e\u0442\u0432\u0435\u0442\u0441\u0442\u0432\u0435\u043d\u043d\u044b\u0439\" -> \"\u041e\u0442\u0432\u0435\u0442\u0441\u0442\u0432\u0435\u043d\u043d\u044b\u0439(-\u0430\u044f)\" ----------------\nconst oldWord2 = \"\u041e\u0442\u0432\u0435\u0442\u0441\u0442\u0432\u0435\u043d\u043d\u044b\u0439 \u2013\";\nconst newWord2 = \"\u041e\u0442\u0432\u0435\u0442\u0441\u0442\u0432\u0435\u043d\u043d\u044b\u0439(-\u0430\u044f) \u2013\";\n\nconst hits2 = body.search(oldWord2, { matchCase: true, matchWholeWord: false });\nhits2.load(\"text\");\nawait context.sync();\n\nif (hits2.items.length > 0) {\n  hits2.items[0].insertText(newWord2, \"Replace\");\n} else {\n  throw new Error(\"Could not find target phrase for change 2\");\n}\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n#\n# Net semantic changes made by the diff (proofErr spell/gram-check markers\n# are cosmetic artifacts of Word's background proofing pass and are not\n# reachable from the Word object model, so we focus on the three actual\n# wording changes the diff introduces):\n#\n#  1. \"...\u0432\u043e \u0432\u0440\u0435\u043c\u044f \u043f\u0440\u043e\u0432\u0435\u0434\u0435\u043d\u0438\u044f \u0441\u043e\u0431\u0440\u0430\u043d\u0438\u0439.\" ->\n#     \"...\u0432\u043e \u0432\u0440\u0435\u043c\u044f \u043f\u0440\u043e\u0432\u0435\u0434\u0435\u043d\u0438\u044f \u0441\u043e\u0431\u0440\u0430\u043d\u0438\u044f, \u0430 \u0442\u0430\u043a\u0436\u0435 \u0433\u0430\u0440\u0430\u043d\u0442\u0438\u0440\u0443\u0435\u043c \u0441\u043e\u0431\u043b\u044e\u0434\u0435\u043d\u0438\u0435 \u043c\u0435\u0440\n#     \u0438\u043d\u0434\u0438\u0432\u0438\u0434\u0443\u0430\u043b\u044c\u043d\u043e\u0439 \u0437\u0430\u0449\u0438\u0442\u044b (\u0434\u043e 50 \u0447\u0435\u043b\u043e\u0432\u0435\u043a).\"\n#  2. \"\u041e\u0442\u0432\u0435\u0442\u0441\u0442\u0432\u0435\u043d\u043d\u044b\u0439 \u2013 {{responsable}}.\" -> \"\u041e\u0442\u0432\u0435\u0442\u0441\u0442\u0432\u0435\u043d\u043d\u044b\u0439(-\u0430\u044f) \u2013 {{responsable}}.\"\n\n$d = $word.ActiveDocument\n\n# --- Change 1: extend the \"\u041c\u044b \u043f\u0440\u0438\u043d\u0438\u043c\u0430\u0435\u043c...\" sentence -----------------\n$oldTail1 = \"\u0432\u043e \u0432\u0440\u0435\u043c\u044f \u043f\u0440\u043e\u0432\u0435\u0434\u0435\u043d\u0438\u044f \u0441\u043e\u0431\u0440\u0430\u043d\u0438\u0439.\"\n$newTail1 = \"\u0432\u043e \u0432\u0440\u0435\u043c\u044f \u043f\u0440\u043e\u0432\u0435\u0434\u0435\u043d\u0438\u044f \u0441\u043e\u0431\u0440\u0430\u043d\u0438\u044f, \u0430 \u0442\u0430\u043a\u0436\u0435 \u0433\u0430\u0440\u0430\u043d\u0442\u0438\u0440\u0443\u0435\u043c \u0441\u043e\u0431\u043b\u044e\u0434\u0435\u043d\u0438\u0435 \u043c\u0435\u0440 \u0438\u043d\u0434\u0438\u0432\u0438\u0434\u0443\u0430\u043b\u044c\u043d\u043e\u0439 \u0437\u0430\u0449\u0438\u0442\u044b (\u0434\u043e 50 \u0447\u0435\u043b\u043e\u0432\u0435\u043a).\"\n\n$rng1 = $d.Content\n$found1 = $rng1.Find.Execute($oldTail1, $true, $false, $false, $false, $false, $true, 1, $false, $newTail1, 2)\nif (-not $found1) {\n    throw \"Could not find target sentence for change 1\"\n}\n\n# --- Change 2: \"\u041e\u0442\u0432\u0435\u0442\u0441\u0442\u0432\u0435\u043d\u043d\u044b\u0439\" -> \"\u041e\u0442\u0432\u0435\u0442\u0441\u0442\u0432\u0435\u043d\u043d\u044b\u0439(-\u0430\u044f)\" ----------------\n$oldWord2 = \"\u041e\u0442\u0432\u0435\u0442\u0441\u0442\u0432\u0435\u043d\u043d\u044b\u0439 \u2013\"\n$newWord2 = \"\u041e\u0442\u0432\u0435\u0442\u0441\u0442\u0432\u0435\u043d\u043d\u044b\u0439(-\u0430\u044f) \u2013\"\n\n$rng2 = $d.Content\n$found2 = $rng2.Find.Execute($oldWord2, $true, $false, $false, $false, $false, $true, 1, $false, $newWord2, 2)\nif (-not $found2) {\n    throw \"Could not find target phrase for change 2\"\n}\n"}
